$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the raw input values - these feed the SUM formulas in G4/G5
# and the ratio formula in G6, which will recalculate automatically.
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 6

# Make this sheet active and move the selection to D4 (matches the
# sheetView selection change in the diff: A10 -> D4)
$ws.Activate()
$ws.Range("D4").Select()

$excel.Calculate()
